$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The underlying data was recomputed with new TPM values. Rows whose
# "Target cluster" (column D) was "ECs" (rows 2, 4, 6, 8) are removed,
# keeping only the rows whose Target cluster is "Resolving-Mac"
# (originally rows 3, 5, 7, 9), which shift up to become rows 2-5.
# Delete from the bottom up so row numbers of rows still to delete do
# not shift.
$ws.Rows.Item(8).Delete()
$ws.Rows.Item(6).Delete()
$ws.Rows.Item(4).Delete()
$ws.Rows.Item(2).Delete()

# Row 2 (was row 3): ECs / Cd80 / Cd28 / Resolving-Mac
$ws.Range("G2").Value = 1.152905666666667
$ws.Range("H2").Value = 3.458717
$ws.Range("I2").Value = 0.06522949989114324
$ws.Range("J2").Value = 0.06522949989114325
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 5.273410666666667
$ws.Range("N2").Value = 15.820232
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 6.079745040260445
$ws.Range("R2").Value = 54.71770536234401
$ws.Range("S2").Value = 0.06522949989114324
$ws.Range("T2").Value = 0.06522949989114325

# Row 3 (was row 5): FAPs / Cd80 / Cd28 / Resolving-Mac
$ws.Range("G3").Value = 3.981869333333333
$ws.Range("H3").Value = 11.945608
$ws.Range("I3").Value = 0.2252875952949142
$ws.Range("J3").Value = 0.2252875952949143
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 5.273410666666667
$ws.Range("N3").Value = 15.820232
$ws.Range("O3").Value = 1
$ws.Range("P3").Value = 1
$ws.Range("Q3").Value = 20.99803221567289
$ws.Range("R3").Value = 188.982289941056
$ws.Range("S3").Value = 0.2252875952949142
$ws.Range("T3").Value = 0.2252875952949143

# Row 4 (was row 7): MuSCs / Cd80 / Cd28 / Resolving-Mac
$ws.Range("G4").Value = 1.517768666666667
$ws.Range("H4").Value = 4.553306
$ws.Range("I4").Value = 0.08587284626968379
$ws.Range("J4").Value = 0.08587284626968379
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 5.273410666666667
$ws.Range("N4").Value = 15.820232
$ws.Range("O4").Value = 1
$ws.Range("P4").Value = 1
$ws.Range("Q4").Value = 8.003817476332445
$ws.Range("R4").Value = 72.03435728699201
$ws.Range("S4").Value = 0.08587284626968379
$ws.Range("T4").Value = 0.08587284626968379

# Row 5 (was row 9): Resolving-Mac / Cd80 / Cd28 / Resolving-Mac
$ws.Range("G5").Value = 11.02206166666667
$ws.Range("H5").Value = 33.066185
$ws.Range("I5").Value = 0.6236100585442587
$ws.Range("J5").Value = 0.6236100585442588
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 5.273410666666667
$ws.Range("N5").Value = 15.820232
$ws.Range("O5").Value = 1
$ws.Range("P5").Value = 1
$ws.Range("Q5").Value = 58.12385756165779
$ws.Range("R5").Value = 523.1147180549201
$ws.Range("S5").Value = 0.6236100585442587
$ws.Range("T5").Value = 0.6236100585442588
